{"js": "// tdf#112118 fixture tweak: replace the page-break run in the first\n// paragraph with the \"_GoBack\" bookmark (moved from the last, now-empty\n// paragraph), leaving that last paragraph empty.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\n// 1) Remove the page-break run from the first paragraph, but keep the\n//    paragraph mark itself: expand a zero-length range across the whole\n//    paragraph content and delete that range (this removes the <w:r> with\n//    the <w:br w:type=\"page\"/> without merging the paragraph into the next\n//    one).\nconst firstStart = firstParagraph.getRange(\"Start\");\nconst firstEnd = firstParagraph.getRange(\"End\");\nconst firstContent = firstStart.expandTo(firstEnd);\nfirstContent.delete();\nawait context.sync();\n\n// 2) Drop the existing \"_GoBack\" bookmark from the last paragraph while its\n//    name is still unique, before we re-create it elsewhere.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 3) Re-insert the \"_GoBack\" bookmark at the very start of the first\n//    paragraph (where the page break used to be).\nconst newBookmarkRange = firstParagraph.getRange(\"Start\");\nnewBookmarkRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# tdf#112118 fixture tweak: replace the page-break run in the first\n# paragraph with the \"_GoBack\" bookmark (moved from the last, now-empty\n# paragraph), leaving that last paragraph empty.\n$d = $word.ActiveDocument\n\n# 1) Drop the existing \"_GoBack\" bookmark (currently wrapping the start of\n#    the last paragraph) while its name is still unique, before we\n#    re-create it elsewhere.\n$d.Bookmarks.Item(\"_GoBack\").Delete()\n\n# 2) Remove the page-break character from the first paragraph, but leave\n#    the paragraph mark (and paragraph) in place.\n$firstParagraph = $d.Paragraphs.Item(1)\n$breakRange = $firstParagraph.Range.Duplicate()\n$breakRange.Collapse(1)          # wdCollapseStart\n[void]$breakRange.MoveEnd(1, 1)  # wdCharacter, grow end by one character\n$breakRange.Delete()\n\n# 3) Re-insert the \"_GoBack\" bookmark at the very start of the first\n#    paragraph (where the page break used to be).\n$firstParagraphAgain = $d.Paragraphs.Item(1)\n$startRange = $firstParagraphAgain.Range.Duplicate()\n$startRange.Collapse(1)   # wdCollapseStart\n$d.Bookmarks.Add(\"_GoBack\", $startRange)\n"}
